# Auto commit at 2025-10-05  8:27:23.15
# Updates the raw metric figures on the "Metrics" sheet (downstream formulas
# on the "today" sheet that reference them - and the TODAY()-1 cell - pick up
# the new values automatically on recalculation), and moves the remembered
# selection on both sheets.

$wb = $excel.ActiveWorkbook

$wsMetrics = $wb.Worksheets.Item("Metrics")
$wsToday   = $wb.Worksheets.Item("today")

# --- Updated source figures on Metrics ---------------------------------
$wsMetrics.Range("B2").Value  = 27386.510000000002
$wsMetrics.Range("B3").Value  = 22729.309999999998
$wsMetrics.Range("B4").Value  = 8534.51
$wsMetrics.Range("B5").Value  = 1109
$wsMetrics.Range("B6").Value  = 4421408.0399999991
$wsMetrics.Range("B7").Value  = 3734630.4100000006
$wsMetrics.Range("B8").Value  = 1287625.3
$wsMetrics.Range("B9").Value  = 171131
$wsMetrics.Range("B10").Value = 32886745.280000001
$wsMetrics.Range("B11").Value = 31009861.199999999
$wsMetrics.Range("B12").Value = 11569338.73
$wsMetrics.Range("B13").Value = 1268759

# --- Selection bookmarks -------------------------------------------------
# Leave "today" as the active sheet when we're done (it was active before),
# but park the Metrics sheet's own remembered selection at D31 first.
$wsMetrics.Select()
$wsMetrics.Range("D31").Select()

$wsToday.Select()
$wsToday.Range("H13").Select()
